# Apply "contingencies with rene fine" edit to lines_states sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B names shift down: line7/line8 inserted after line6, extr1..extr8
# slide down into rows 10..17 (two new rows appended at the end).
$names = @("line1","line2","line3","line4","line5","line6","line7","line8","extr1","extr2","extr3","extr4","extr5","extr6","extr7","extr8")

# Row data (C, D, E) for rows 2..17, matching the final target state.
$cVals = @(7, 9, 8, 8, 10, 12, 14, 16, 5, 5, 10, 7, 9, 7, 5, 8)
$dVals = @(9, 8, 10, 11, 5, 8, 11, 9, 12, 9, 11, 8, 11, 11, 7, 5)
$eVals = @($true, $true, $true, $false, $true, $true, $true, $true, $true, $true, $false, $true, $true, $true, $true, $true)

# Make sure rows 16 and 17 exist with the same formatting as the other data
# rows (bold/centered/bordered style) before writing into column A.
$ws.Cells.Item(15, 1).Copy() | Out-Null
$ws.Cells.Item(16, 1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(15, 1).Copy() | Out-Null
$ws.Cells.Item(17, 1).PasteSpecial(-4122) | Out-Null

for ($i = 0; $i -lt 16; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 1).Value = $i
    $ws.Cells.Item($row, 2).Value = $names[$i]
    $ws.Cells.Item($row, 3).Value = $cVals[$i]
    $ws.Cells.Item($row, 4).Value = $dVals[$i]
    $ws.Cells.Item($row, 5).Value = $eVals[$i]
}

Write-Host "done"
